# The sheet previously carried a leading, unheaded numeric column (A) that
# duplicated the GENE id already present in the last column. Remove it so
# the real headers (EL_Astral15, FNRATE_PHYLONET, TAXON, MODEL_CONDITION,
# GENE) line up with columns A:E, then fix the typo'd header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray leading column; everything else (B:F) shifts left to A:E.
$ws.Columns("A").Delete()

# Header row: "MODEL_CONDITION" -> "MODELCONDITION" (now in column D).
$ws.Range("D1").Value = "MODELCONDITION"
